# Update gh-pages output (合肥-漫展信息.xlsx) to data generated at 456a3b4.
# The "合肥·漫有引力动漫游戏嘉年华" (2024-09-21) listing has fallen off the
# bottom of the feed; every later row shifts up by one and the data that
# used to be in the last row disappears (the two affected sheets each lose
# exactly one row at the bottom).

$wb = $excel.ActiveWorkbook

function Set-EventRow($ws, $Row, $B, $C, $D, $E, $F, $G, $H, $I) {
    # Column B holds plain "yyyy-MM-dd" text in the source data, not a real
    # date value, so force text formatting first to stop Excel from
    # auto-converting it to a date serial number.
    $ws.Range("B$Row").NumberFormat = "@"
    $ws.Range("B$Row").Value = $B
    $ws.Range("C$Row").Value = $C
    $ws.Range("D$Row").Value = $D
    $ws.Range("E$Row").Value = $E
    $ws.Range("F$Row").Value = $F
    $ws.Range("G$Row").Value = $G
    $ws.Range("H$Row").Value = $H
    $ws.Range("I$Row").Value = $I
}

# ---------------------------------------------------------------------------
# Sheet "展览" (exhibitions): rows 2-13 -> rows 2-12 (row 13 disappears)
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

Set-EventRow $wsExpo 2  "2024-09-15" "蜀山·银泰百货高新店-2024漫趣地带嘉年华（免费）" "高新区望江西路888号 银泰百货（高新店）" "2024.09.15 10:00-10.02 22:00" 276 30 "https://show.bilibili.com/platform/detail.html?id=91869" "//i2.hdslb.com/bfs/openplatform/202409/JDGIWMyd1725422862878.png"
Set-EventRow $wsExpo 3  "2024-10-01" "合肥·星域动漫游戏嘉年华" "新站区东方大道288号 少荃体育中心" "2024.10.01 10:00-10.01 17:00" 15 58 "https://show.bilibili.com/platform/detail.html?id=91878" "//i0.hdslb.com/bfs/openplatform/202409/NOg6Wwjh1725121441581.png"
Set-EventRow $wsExpo 4  "2024-10-01" "合肥·第十五届次元之门动漫游戏博览会" "南京路与庐州大道交汇处 合肥滨湖国际会展中心" "2024.10.01 09:30-10.02 17:30" 7679 70 "https://show.bilibili.com/platform/detail.html?id=91133" "//i1.hdslb.com/bfs/openplatform/202408/PlcqtYWR1724315434068.jpeg"
Set-EventRow $wsExpo 5  "2024-10-01" "合肥·首届AT次元时代动漫游戏嘉年华" "凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心" "2024.10.01 09:30-10.03 17:00" 5589 68 "https://show.bilibili.com/platform/detail.html?id=90908" "//i2.hdslb.com/bfs/openplatform/202409/RqtCRIaH1726800618582.jpeg"
Set-EventRow $wsExpo 6  "2024-10-04" "合肥·Holic动漫游戏展" "庐州大道800号 合肥融创茂" "2024.10.04 10:00-10.06 17:00" 461 55 "https://show.bilibili.com/platform/detail.html?id=92061" "//i1.hdslb.com/bfs/openplatform/202409/AZ0LsUce1725522015668.jpeg"
Set-EventRow $wsExpo 7  "2024-10-04" "合肥·乐帮•崩原铁绝only同人首展" "丹霞路488号金星商业城三楼 迷鹿轰趴" "2024.10.04 10:00-10.05 16:30" 71 58 "https://show.bilibili.com/platform/detail.html?id=91524" "//i2.hdslb.com/bfs/openplatform/202408/739I7YRr1724912450704.png"
Set-EventRow $wsExpo 8  "2024-10-06" "合肥·星月动漫游戏展" "灵石路与皇藏峪路交叉口西南10米安徽百事兴电气有限公司院内2栋厂房2层 兄弟篮球俱乐部" "2024.10.06 10:00-10.06 17:00" 10 29.9 "https://show.bilibili.com/platform/detail.html?id=91958" "//i2.hdslb.com/bfs/openplatform/202409/mgB8U6bN1725361649767.jpeg"
Set-EventRow $wsExpo 9  "2024-10-06" "合肥·首届火影忍者同人only" "长江东路金太阳家具广场南门二楼 优极篮球馆" "2024.10.06 09:30-10.06 17:30" 66 75 "https://show.bilibili.com/platform/detail.html?id=91658" "//i0.hdslb.com/bfs/openplatform/202408/f8ylbskH1725027552569.jpeg"
Set-EventRow $wsExpo 10 "2024-10-26" "合肥·W·A第五人格同人only2.0" "莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)" "2024.10.26 09:30-10.26 17:00" 249 68 "https://show.bilibili.com/platform/detail.html?id=91123" "//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png"
Set-EventRow $wsExpo 11 "2024-11-16" "合肥·第九届环形宇宙动漫游戏嘉年华" "南京路与庐州大道交汇处 合肥滨湖国际会展中心" "2024.11.16 09:30-11.17 17:00" 209 70 "https://show.bilibili.com/platform/detail.html?id=92565" "//i1.hdslb.com/bfs/openplatform/202409/WuHNviSs1726646410055.jpeg"
Set-EventRow $wsExpo 12 "2024-11-17" "合肥·MAX特摄同人only2.0" "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间" "2024.11.17 10:00-11.17 18:00" 52 60 "https://show.bilibili.com/platform/detail.html?id=92498" "//i0.hdslb.com/bfs/openplatform/202409/R4WJxlQe1726230330813.jpeg"

# Drop the now-surplus trailing row (old row 13), shrinking the used range.
$wsExpo.Rows.Item(13).Delete()

# ---------------------------------------------------------------------------
# Sheet "全部类型" (all types): rows 2-16 -> rows 2-15 (row 16 disappears)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

Set-EventRow $wsAll 2  "2024-09-15" "蜀山·银泰百货高新店-2024漫趣地带嘉年华（免费）" "高新区望江西路888号 银泰百货（高新店）" "2024.09.15 10:00-10.02 22:00" 276 30 "https://show.bilibili.com/platform/detail.html?id=91869" "//i2.hdslb.com/bfs/openplatform/202409/JDGIWMyd1725422862878.png"
Set-EventRow $wsAll 3  "2024-10-01" "合肥·星域动漫游戏嘉年华" "新站区东方大道288号 少荃体育中心" "2024.10.01 10:00-10.01 17:00" 15 58 "https://show.bilibili.com/platform/detail.html?id=91878" "//i0.hdslb.com/bfs/openplatform/202409/NOg6Wwjh1725121441581.png"
Set-EventRow $wsAll 4  "2024-10-01" "合肥·第十五届次元之门动漫游戏博览会" "南京路与庐州大道交汇处 合肥滨湖国际会展中心" "2024.10.01 09:30-10.02 17:30" 7679 70 "https://show.bilibili.com/platform/detail.html?id=91133" "//i1.hdslb.com/bfs/openplatform/202408/PlcqtYWR1724315434068.jpeg"
Set-EventRow $wsAll 5  "2024-10-01" "合肥·首届AT次元时代动漫游戏嘉年华" "凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心" "2024.10.01 09:30-10.03 17:00" 5589 68 "https://show.bilibili.com/platform/detail.html?id=90908" "//i2.hdslb.com/bfs/openplatform/202409/RqtCRIaH1726800618582.jpeg"
Set-EventRow $wsAll 6  "2024-10-04" "合肥·Holic动漫游戏展" "庐州大道800号 合肥融创茂" "2024.10.04 10:00-10.06 17:00" 461 55 "https://show.bilibili.com/platform/detail.html?id=92061" "//i1.hdslb.com/bfs/openplatform/202409/AZ0LsUce1725522015668.jpeg"
Set-EventRow $wsAll 7  "2024-10-04" "合肥·乐帮•崩原铁绝only同人首展" "丹霞路488号金星商业城三楼 迷鹿轰趴" "2024.10.04 10:00-10.05 16:30" 71 58 "https://show.bilibili.com/platform/detail.html?id=91524" "//i2.hdslb.com/bfs/openplatform/202408/739I7YRr1724912450704.png"
Set-EventRow $wsAll 8  "2024-10-06" "合肥·星月动漫游戏展" "灵石路与皇藏峪路交叉口西南10米安徽百事兴电气有限公司院内2栋厂房2层 兄弟篮球俱乐部" "2024.10.06 10:00-10.06 17:00" 10 29.9 "https://show.bilibili.com/platform/detail.html?id=91958" "//i2.hdslb.com/bfs/openplatform/202409/mgB8U6bN1725361649767.jpeg"
Set-EventRow $wsAll 9  "2024-10-06" "合肥·首届火影忍者同人only" "长江东路金太阳家具广场南门二楼 优极篮球馆" "2024.10.06 09:30-10.06 17:30" 66 75 "https://show.bilibili.com/platform/detail.html?id=91658" "//i0.hdslb.com/bfs/openplatform/202408/f8ylbskH1725027552569.jpeg"
Set-EventRow $wsAll 10 "2024-10-26" "合肥·W·A第五人格同人only2.0" "莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)" "2024.10.26 09:30-10.26 17:00" 249 68 "https://show.bilibili.com/platform/detail.html?id=91123" "//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png"
Set-EventRow $wsAll 11 "2024-10-26" "合肥·《四月是你的谎言》—“公生”与“薰”的钢琴小提琴唯美经典音乐集" "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院" "2024.10.26 19:30-10.26 21:00" 83 80 "https://show.bilibili.com/platform/detail.html?id=90322" "//i2.hdslb.com/bfs/openplatform/202408/BiVgXUKH1722824304648.jpeg"
Set-EventRow $wsAll 12 "2024-11-09" "合肥·一生必听的钢琴曲—“从巴赫 · 莫扎特到肖邦 · 李斯特”钢琴圣手谭小棠独奏音乐会" "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院" "2024.11.09 19:30-11.09 21:00" 5 64 "https://show.bilibili.com/platform/detail.html?id=90593" "//i2.hdslb.com/bfs/openplatform/202408/SYfLxnO21723442234232.jpeg"
Set-EventRow $wsAll 13 "2024-11-16" "合肥·第九届环形宇宙动漫游戏嘉年华" "南京路与庐州大道交汇处 合肥滨湖国际会展中心" "2024.11.16 09:30-11.17 17:00" 209 70 "https://show.bilibili.com/platform/detail.html?id=92565" "//i1.hdslb.com/bfs/openplatform/202409/WuHNviSs1726646410055.jpeg"
Set-EventRow $wsAll 14 "2024-11-17" "合肥·MAX特摄同人only2.0" "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间" "2024.11.17 10:00-11.17 18:00" 52 60 "https://show.bilibili.com/platform/detail.html?id=92498" "//i0.hdslb.com/bfs/openplatform/202409/R4WJxlQe1726230330813.jpeg"
Set-EventRow $wsAll 15 "2024-12-07" "合肥·一生必听的古典系列《钟》—超技钢琴曲炫彩音乐会" "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院" "2024.12.07 19:30-12.07 21:00" 2 64 "https://show.bilibili.com/platform/detail.html?id=91608" "//i0.hdslb.com/bfs/openplatform/202408/wiLiWoeM1725005636569.jpeg"

# Drop the now-surplus trailing row (old row 16), shrinking the used range.
$wsAll.Rows.Item(16).Delete()
